# Add the new "Anuales" worksheet after the existing "Sheet3" sheet
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Anuales"

# Match the page margins used elsewhere in the workbook (in points: 1 inch = 72pt)
$pageSetup = $newSheet.PageSetup
$pageSetup.LeftMargin = 54
$pageSetup.RightMargin = 54
$pageSetup.TopMargin = 72
$pageSetup.BottomMargin = 72
$pageSetup.HeaderMargin = 36
$pageSetup.FooterMargin = 36

# Data block (40 rows x 9 cols) pulled from the workbook history
$csvData = @"
0,5.0000000000000001E-4,5.0000000000000001E-4,5.0000000000000001E-4,5.0000000000000001E-4,1.9084364924755301E-3,1.97801325374013,0.319279253862434,9
0,5.0000000000000001E-4,5.0000000000000001E-4,5.0000000000000001E-4,5.0000000000000001E-4,4.14399041389197E-4,1.3842632546911899,8.8271214249260904E-2,3
0,1.5E-3,1.5E-3,1.5E-3,1.5E-3,5.6756361826644898E-3,1.9772372100951101,0.366250800748393,8
0,1.5E-3,1.5E-3,1.5E-3,1.5E-3,1.2378691610517099E-3,1.38562017648697,2.29904917085277E-2,1
0,2.5000000000000001E-3,2.5000000000000001E-3,2.5000000000000001E-3,2.5000000000000001E-3,9.5079668990434903E-3,1.9731193499691599,0.20442020318512,10
0,2.5000000000000001E-3,2.5000000000000001E-3,2.5000000000000001E-3,2.5000000000000001E-3,2.03156968801929E-3,1.3928722702365399,5.5151382390454001E-2,1
0,5.0000000000000001E-3,5.0000000000000001E-3,5.0000000000000001E-3,5.0000000000000001E-3,1.89463919050324E-2,2.0431532253482301,0.23035020519579399,6
0,5.0000000000000001E-3,5.0000000000000001E-3,5.0000000000000001E-3,5.0000000000000001E-3,3.9856211393866801E-3,1.3298664082140901,8.7411678396951201E-3,1
0,1E-3,0,0,0,1.69583340781817E-3,1.2345137935638399,3.9549073674308703E-2,2
0,1E-3,0,0,0,2.5765840575334701E-4,0.96475386948401798,7.7338772307235501E-2,1
0,3.0000000000000001E-3,0,0,0,4.3633986287244103E-3,1.23780460980101,7.6267822984057404E-2,4
0,3.0000000000000001E-3,0,0,0,6.7628329423236099E-4,0.940123483915974,0.12536843791367899,3
0,5.0000000000000001E-3,0,0,0,8.1280100393522695E-3,1.2842206120036499,0.42925004171654102,9
0,5.0000000000000001E-3,0,0,0,1.24963585769882E-3,0.98732501405636897,0.20326390318162499,6
0,0.01,0,0,0,1.5971900222063899E-2,1.2855219107692899,0.18817094875634399,7
0,0.01,0,0,0,2.4719067362624502E-3,0.95553050297585196,5.8433574616541299E-2,1
0,0,1E-3,0,0,1.3555548185636001E-3,1.2552323978250599,0.13811872675374001,7
0,0,1E-3,0,0,3.6505803855203502E-4,1.09157700417993,0.495432386846555,10
0,0,3.0000000000000001E-3,0,0,4.12467851868413E-3,1.28392029356869,0.20021213639145299,8
0,0,3.0000000000000001E-3,0,0,1.1015811990259599E-3,1.1354257523743101,0.380051944209082,10
0,0,5.0000000000000001E-3,0,0,6.5513833329157601E-3,1.1509094070683099,0.482321477719867,10
0,0,5.0000000000000001E-3,0,0,1.74711534527732E-3,1.0361747532329399,0.72770226072784405,10
0,0,0.01,0,0,1.24209187495645E-2,1.1541186370910399,0.366928762999845,8
0,0,0.01,0,0,3.3132711200353601E-3,1.0316779911865599,0.471135948194289,9
0,0,0,1E-3,0,2.2517644538677202E-3,1.81194000118803,0.377292951717138,10
0,0,0,1E-3,0,4.5023484835890801E-4,1.29078503723055,7.6069782355433196E-3,0
0,0,0,3.0000000000000001E-3,0,6.4459294822511396E-3,1.7891098062373201,0.11681436023526801,6
0,0,0,3.0000000000000001E-3,0,1.2960439289970301E-3,1.29906901631668,5.3023789899728502E-2,3
0,0,0,5.0000000000000001E-3,0,1.0761387644495799E-2,1.82876966438887,0.184027563233399,6
0,0,0,5.0000000000000001E-3,0,2.1809825231056401E-3,1.33682208515547,0.13084481263937101,4
0,0,0,0.01,0,2.1233137800906798E-2,1.8227401471585201,0.26109540051616298,7
0,0,0,0.01,0,4.4626330059494897E-3,1.28362101218998,0.15128428418746501,5
0,0,0,0,1E-3,1.6116871624178399E-3,1.3119586405747801,0.13780480254199401,8
0,0,0,0,1E-3,3.3237147030305901E-4,0.99347207762484702,0.37435124325101798,8
0,0,0,0,3.0000000000000001E-3,4.9606787971784704E-3,1.2080574104197199,0.28646741346771398,8
0,0,0,0,3.0000000000000001E-3,1.0233365530762E-3,0.92773271840362403,0.465696365531371,10
0,0,0,0,5.0000000000000001E-3,8.5055842915500099E-3,1.3039083327249099,0.19240098971160299,6
0,0,0,0,5.0000000000000001E-3,1.7711685868725801E-3,0.97716675006880604,9.3675515967147002E-2,4
0,0,0,0,0.01,1.6804709817280598E-2,1.2639282717301601,0.18995161231087199,5
0,0,0,0,0.01,3.5202625921716799E-3,0.95105048190653696,0.454305213235139,8
"@

$lines = $csvData -split "`n"
$nRows = $lines.Count
$nCols = 9
$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    $fields = $lines[$i].Trim() -split ","
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i,$j] = [double]$fields[$j]
    }
}

$targetRange = $newSheet.Range($newSheet.Cells.Item(1,1), $newSheet.Cells.Item($nRows, $nCols))
$targetRange.Value = $arr

# Make "Anuales" the active sheet/tab and set its selection to K8
$newSheet.Activate()
$newSheet.Range("K8").Select()
